$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "68+6="
$t.Cell(1, 2).Range.Text = "92-6="
$t.Cell(1, 3).Range.Text = "87-68="
$t.Cell(1, 4).Range.Text = "92-44="
$t.Cell(1, 5).Range.Text = "39+33="
$t.Cell(2, 1).Range.Text = "45+7="
$t.Cell(2, 2).Range.Text = "75+6="
$t.Cell(2, 3).Range.Text = "81-79="
$t.Cell(2, 4).Range.Text = "66-39="
$t.Cell(2, 5).Range.Text = "64-26="
$t.Cell(3, 1).Range.Text = "39+19="
$t.Cell(3, 2).Range.Text = "62-36="
$t.Cell(3, 3).Range.Text = "81-33="
$t.Cell(3, 4).Range.Text = "66-58="
$t.Cell(3, 5).Range.Text = "65-19="
$t.Cell(4, 1).Range.Text = "93-19="
$t.Cell(4, 2).Range.Text = "68+6="
$t.Cell(4, 3).Range.Text = "18+45="
$t.Cell(4, 4).Range.Text = "73-6="
$t.Cell(4, 5).Range.Text = "90-22="
$t.Cell(5, 1).Range.Text = "92-23="
$t.Cell(5, 2).Range.Text = "18+67="
$t.Cell(5, 3).Range.Text = "74-29="
$t.Cell(5, 4).Range.Text = "19+45="
$t.Cell(5, 5).Range.Text = "46+7="
$t.Cell(6, 1).Range.Text = "80-67="
$t.Cell(6, 2).Range.Text = "53+29="
$t.Cell(6, 3).Range.Text = "19+19="
$t.Cell(6, 4).Range.Text = "41-2="
$t.Cell(6, 5).Range.Text = "84-25="
$t.Cell(7, 1).Range.Text = "91-32="
$t.Cell(7, 2).Range.Text = "75-69="
$t.Cell(7, 3).Range.Text = "39+24="
$t.Cell(7, 4).Range.Text = "25-6="
$t.Cell(7, 5).Range.Text = "9+38="
$t.Cell(8, 1).Range.Text = "70-52="
$t.Cell(8, 2).Range.Text = "97-78="
$t.Cell(8, 3).Range.Text = "90-3="
$t.Cell(8, 4).Range.Text = "27+29="
$t.Cell(8, 5).Range.Text = "70-45="
$t.Cell(9, 1).Range.Text = "77+18="
$t.Cell(9, 2).Range.Text = "71-29="
$t.Cell(9, 3).Range.Text = "50-33="
$t.Cell(9, 4).Range.Text = "16-9="
$t.Cell(9, 5).Range.Text = "18+55="
$t.Cell(10, 1).Range.Text = "95-16="
$t.Cell(10, 2).Range.Text = "24+18="
$t.Cell(10, 3).Range.Text = "42-5="
$t.Cell(10, 4).Range.Text = "49+24="
$t.Cell(10, 5).Range.Text = "87-69="
$t.Cell(11, 1).Range.Text = "89+6="
$t.Cell(11, 2).Range.Text = "32-27="
$t.Cell(11, 3).Range.Text = "55-39="
$t.Cell(11, 4).Range.Text = "7+29="
$t.Cell(11, 5).Range.Text = "49+29="
$t.Cell(12, 1).Range.Text = "76-39="
$t.Cell(12, 2).Range.Text = "9+4="
$t.Cell(12, 3).Range.Text = "33-24="
$t.Cell(12, 4).Range.Text = "6+49="
$t.Cell(12, 5).Range.Text = "33-7="
$t.Cell(13, 1).Range.Text = "77+4="
$t.Cell(13, 2).Range.Text = "75-67="
$t.Cell(13, 3).Range.Text = "6+87="
$t.Cell(13, 4).Range.Text = "49+36="
$t.Cell(13, 5).Range.Text = "90-53="
$t.Cell(14, 1).Range.Text = "17+5="
$t.Cell(14, 2).Range.Text = "42-23="
$t.Cell(14, 3).Range.Text = "92-8="
$t.Cell(14, 4).Range.Text = "16+19="
$t.Cell(14, 5).Range.Text = "93-36="
$t.Cell(15, 1).Range.Text = "81-24="
$t.Cell(15, 2).Range.Text = "82-14="
$t.Cell(15, 3).Range.Text = "29+35="
$t.Cell(15, 4).Range.Text = "42+49="
$t.Cell(15, 5).Range.Text = "71-3="
$t.Cell(16, 1).Range.Text = "9+72="
$t.Cell(16, 2).Range.Text = "60-19="
$t.Cell(16, 3).Range.Text = "81-72="
$t.Cell(16, 4).Range.Text = "19+52="
$t.Cell(16, 5).Range.Text = "92-55="
$t.Cell(17, 1).Range.Text = "27+9="
$t.Cell(17, 2).Range.Text = "62+19="
$t.Cell(17, 3).Range.Text = "6+69="
$t.Cell(17, 4).Range.Text = "40-31="
$t.Cell(17, 5).Range.Text = "28+65="
$t.Cell(18, 1).Range.Text = "52-48="
$t.Cell(18, 2).Range.Text = "9+65="
$t.Cell(18, 3).Range.Text = "95-36="
$t.Cell(18, 4).Range.Text = "30-27="
$t.Cell(18, 5).Range.Text = "2+79="
$t.Cell(19, 1).Range.Text = "69+12="
$t.Cell(19, 2).Range.Text = "7+26="
$t.Cell(19, 3).Range.Text = "93-7="
$t.Cell(19, 4).Range.Text = "41-37="
$t.Cell(19, 5).Range.Text = "26+6="
$t.Cell(20, 1).Range.Text = "15+7="
$t.Cell(20, 2).Range.Text = "84-78="
$t.Cell(20, 3).Range.Text = "96-29="
$t.Cell(20, 4).Range.Text = "84-46="
$t.Cell(20, 5).Range.Text = "37-29="
